$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 4.576389969013007
    "C2" = 0.3342714604665957
    "D2" = 0.006836383813318747
    "E2" = 0.05039595700303684
    "F2" = 4.735664047197844
    "J2" = 0.1416895320638467
    "L2" = 0.4414515150752294
    "N2" = 2.757242386420714
    "B3" = 4.455837958332722
    "C3" = 0.3080797494909575
    "D3" = 0.006087666368689071
    "E3" = 0.05047652714794171
    "F3" = 4.707275549390204
    "J3" = 0.142404487373808
    "L3" = 0.4365187045586794
    "N3" = 2.774832669723523
    "B4" = 4.3845890641295
    "C4" = 0.2922051183752785
    "D4" = 0.005626340931215168
    "E4" = 0.05053027515297426
    "F4" = 4.69224851944611
    "J4" = 0.1428772730636485
    "L4" = 0.4337185642581858
    "N4" = 2.786423068807878
    "B5" = 4.356248667022726
    "C5" = 0.2857874830917808
    "D5" = 0.005437879225244302
    "E5" = 0.05055325506186159
    "F5" = 4.686727119348774
    "J5" = 0.1430784426896725
    "L5" = 0.4326349052858802
    "N5" = 2.791344288055754
    "B6" = 4.351584614419039
    "C6" = 0.2847249244072145
    "D6" = 0.005406555212061903
    "E6" = 0.05055713596214928
    "F6" = 4.685846605860192
    "J6" = 0.1431123606589644
    "L6" = 0.4324584307085217
    "N6" = 2.792173396664921
    "B7" = 4.384204048535025
    "C7" = 0.2921183606683257
    "D7" = 0.005623801248152205
    "E7" = 0.0505305807044486
    "F7" = 4.692171620306439
    "J7" = 0.1428799516596566
    "L7" = 0.4337037172449385
    "N7" = 2.786488637009114
    "B8" = 4.534246999201685
    "C8" = 0.3251971288768232
    "D8" = 0.006578521943545468
    "E8" = 0.05042285066558033
    "F8" = 4.725375482668454
    "J8" = 0.1419290371532487
    "L8" = 0.4397031855137783
    "N8" = 2.763143219145007
    "B9" = 4.850603414006628
    "C9" = 0.3917435685590362
    "D9" = 0.00844087024803386
    "E9" = 0.05024546993411666
    "F9" = 4.809668012593818
    "J9" = 0.1403322143123873
    "L9" = 0.4532866387674801
    "N9" = 2.723656710720775
    "B10" = 5.096745335844389
    "C10" = 0.4417131830362564
    "D10" = 0.009807277632475575
    "E10" = 0.05013572106777198
    "F10" = 4.883454828677458
    "J10" = 0.1393219986152801
    "L10" = 0.4643832081064403
    "N10" = 2.698519053072658
    "B11" = 5.211752386229819
    "C11" = 0.464691975589858
    "D11" = 0.01042939837231671
    "E11" = 0.05009024535388051
    "F11" = 4.919634178634141
    "J11" = 0.1388977419361481
    "L11" = 0.4696758015430618
    "N11" = 2.687932273078744
    "B12" = 5.255742883157268
    "C12" = 0.4734299227253587
    "D12" = 0.01066512940509767
    "E12" = 0.05007366363443322
    "F12" = 4.93371301120888
    "J12" = 0.1387421577391592
    "L12" = 0.4717152914154212
    "N12" = 2.684046041560734
    "B13" = 5.246249128182171
    "C13" = 0.471546418573098
    "D13" = 0.01061435282173306
    "E13" = 0.05007720640074709
    "F13" = 4.930664005281812
    "J13" = 0.1387754400170955
    "L13" = 0.4712744789963921
    "N13" = 2.684877539615968
    "B14" = 5.215362677103144
    "C14" = 0.4654101180308317
    "D14" = 0.01044878874565569
    "E14" = 0.05008886836530735
    "F14" = 4.920784851331149
    "J14" = 0.1388848402837546
    "L14" = 0.4698428834909976
    "N14" = 2.687610086029721
    "B15" = 5.196501210637393
    "C15" = 0.4616562183531983
    "D15" = 0.01034739727096223
    "E15" = 0.05009609483968502
    "F15" = 4.914782957005684
    "J15" = 0.1389525116278669
    "L15" = 0.4689705903443979
    "N15" = 2.689299857872101
    "B16" = 5.089290721179623
    "C16" = 0.4402165110670921
    "D16" = 0.009766637349486018
    "E16" = 0.05013878247095127
    "F16" = 4.881143240022737
    "J16" = 0.1393504347728687
    "L16" = 0.4640422581818058
    "N16" = 2.699228060650043
    "B17" = 5.024300410456988
    "C17" = 0.4271278696602963
    "D17" = 0.009410546862810065
    "E17" = 0.05016610893304918
    "F17" = 4.861177609923999
    "J17" = 0.13960358564686
    "L17" = 0.461081638791768
    "N17" = 2.705536529286334
    "B18" = 4.987205250342186
    "C18" = 0.4196228631919894
    "D18" = 0.009205779137921866
    "E18" = 0.05018224525729997
    "F18" = 4.849939724228818
    "J18" = 0.139752513669146
    "L18" = 0.4594017949643074
    "N18" = 2.709244810053477
    "B19" = 4.974694411488599
    "C19" = 0.4170857630961109
    "D19" = 0.00913645409225694
    "E19" = 0.0501877807087977
    "F19" = 4.846176907089301
    "J19" = 0.1398035088443077
    "L19" = 0.4588369801846426
    "N19" = 2.710514056278882
    "B20" = 5.031189163502233
    "C20" = 0.4285187674723829
    "D20" = 0.009448448024087952
    "E20" = 0.05016315664032023
    "F20" = 4.863277524870483
    "J20" = 0.1395762934766758
    "L20" = 0.4613944178772726
    "N20" = 2.704856715986608
    "B21" = 5.224422813265846
    "C21" = 0.4672115044746192
    "D21" = 0.01049741436190743
    "E21" = 0.05008542562953699
    "F21" = 4.923676304184596
    "J21" = 0.1388525691258558
    "L21" = 0.4702624188607416
    "N21" = 2.68680413345335
    "B22" = 5.353277385893421
    "C22" = 0.4927118644701522
    "D22" = 0.01118386706581731
    "E22" = 0.05003834782892291
    "F22" = 4.965357868210305
    "J22" = 0.1384091392832261
    "L22" = 0.47626397875419
    "N22" = 2.675721575942546
    "B23" = 5.284269452310923
    "C23" = 0.4790821352793841
    "D23" = 0.01081738915263486
    "E23" = 0.0500631336758437
    "F23" = 4.942908742006466
    "J23" = 0.1386431017301319
    "L23" = 0.4730419636801457
    "N23" = 2.68157079608072
    "B24" = 5.028073921757596
    "C24" = 0.4278898808727831
    "D24" = 0.009431313059273805
    "E24" = 0.05016449004498047
    "F24" = 4.8623274040458
    "J24" = 0.1395886217131057
    "L24" = 0.4612529411390653
    "N24" = 2.70516380591198
    "B25" = 4.762628969091622
    "C25" = 0.3735556747829492
    "D25" = 0.007937684988863936
    "E25" = 0.05028983783153101
    "F25" = 4.784794499094176
    "J25" = 0.1407355528816367
    "L25" = 0.4494164429841589
    "N25" = 2.733661364364977
}

foreach ($cellRef in $values.Keys) {
    $ws.Range($cellRef).Value = $values[$cellRef]
}

Write-Host "Updated $($values.Count) cells"
